$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3, 17, 11, "2025-03-27 08:22:27", 41500.0415, 0, "O122"),
    @(4, 17, 10, "2025-03-27 15:28:01", 59000.059,  0, "O827"),
    @(5, 18, 11, "2025-03-27 17:20:34", 7500.075,   0, "O1020")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
